# Re-point heading style formatting one level down (Heading2 gets the old
# Heading3 size, Heading3 gets the old Heading4 size, etc.), drop the
# themeShade on Heading1's color so every heading uses the same accent1
# blue, make Heading4 italic instead of bold, and make Heading5 plain
# (no italic) instead of italic.

$d = $word.ActiveDocument

$styles = $d.Styles

$h1 = $styles.Item("Heading 1")
$h2 = $styles.Item("Heading 2")
$h3 = $styles.Item("Heading 3")
$h4 = $styles.Item("Heading 4")
$h5 = $styles.Item("Heading 5")

# Heading 1: same color, but without the themeShade-darkened variant -
# 0x4F81BD (accent1), same RGB value the other headings already use.
# Word COM colors are packed 0x00BBGGRR, so build it from the R/G/B bytes
# rather than writing the hex string "as-is".
$h1.Font.Color = 0x4F + (0x81 * 0x100) + (0xBD * 0x10000)

# Heading 2: 16pt -> 14pt.
$h2.Font.Size = 14
$h2.Font.SizeBi = 14

# Heading 3: 14pt -> 12pt.
$h3.Font.Size = 12
$h3.Font.SizeBi = 12

# Heading 4: bold -> italic (still 12pt).
$h4.Font.Bold = $false
$h4.Font.Italic = $true

# Heading 5: no longer italic.
$h5.Font.Italic = $false
